$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new "Save" header in column H, matching the formatting of the
# existing header cells (e.g. G1: bold font, border, centered).
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add the corresponding numeric value in row 2.
$ws.Range("H2").Value = 1
